$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = "sd"
$ws.Range("J5").Value = "Statement-non-opinion"
$ws.Range("I22").Value = "sd"
$ws.Range("J22").Value = "Statement-non-opinion"
$ws.Range("I25").Value = "aa"
$ws.Range("J25").Value = "Agree/Accept"
$ws.Range("I40").Value = "b"
$ws.Range("J40").Value = "Acknowledge (Backchannel)"
$ws.Range("I51").Value = "%"
$ws.Range("J51").Value = "Uninterpretable"
$ws.Range("I52").Value = "%"
$ws.Range("J52").Value = "Uninterpretable"
$ws.Range("I59").Value = "ba"
$ws.Range("J59").Value = "Appreciation"
$ws.Range("I94").Value = "%"
$ws.Range("J94").Value = "Uninterpretable"
$ws.Range("I108").Value = "sd"
$ws.Range("J108").Value = "Statement-non-opinion"
$ws.Range("I109").Value = "ba"
$ws.Range("J109").Value = "Appreciation"
$ws.Range("I120").Value = "sd"
$ws.Range("J120").Value = "Statement-non-opinion"
$ws.Range("I132").Value = "sd"
$ws.Range("J132").Value = "Statement-non-opinion"
$ws.Range("I133").Value = "sd"
$ws.Range("J133").Value = "Statement-non-opinion"
$ws.Range("I166").Value = "ba"
$ws.Range("J166").Value = "Appreciation"
$ws.Range("I173").Value = "sd"
$ws.Range("J173").Value = "Statement-non-opinion"
$ws.Range("I175").Value = "ba"
$ws.Range("J175").Value = "Appreciation"
$ws.Range("I181").Value = "aa"
$ws.Range("J181").Value = "Agree/Accept"
$ws.Range("I189").Value = "sd"
$ws.Range("J189").Value = "Statement-non-opinion"
$ws.Range("I193").Value = "sd"
$ws.Range("J193").Value = "Statement-non-opinion"
$ws.Range("I194").Value = "ba"
$ws.Range("J194").Value = "Appreciation"
$ws.Range("I196").Value = "sv"
$ws.Range("J196").Value = "Statement-opinion"
$ws.Range("I197").Value = "aa"
$ws.Range("J197").Value = "Agree/Accept"
$ws.Range("I214").Value = "sd"
$ws.Range("J214").Value = "Statement-non-opinion"
$ws.Range("I220").Value = "sv"
$ws.Range("J220").Value = "Statement-opinion"
$ws.Range("I224").Value = "sd"
$ws.Range("J224").Value = "Statement-non-opinion"
$ws.Range("I239").Value = "sv"
$ws.Range("J239").Value = "Statement-opinion"
$ws.Range("I243").Value = "sv"
$ws.Range("J243").Value = "Statement-opinion"
$ws.Range("I246").Value = "b"
$ws.Range("J246").Value = "Acknowledge (Backchannel)"
$ws.Range("I247").Value = "sv"
$ws.Range("J247").Value = "Statement-opinion"
$ws.Range("I249").Value = "aa"
$ws.Range("J249").Value = "Agree/Accept"
$ws.Range("I254").Value = "%"
$ws.Range("J254").Value = "Uninterpretable"
$ws.Range("I259").Value = "sd"
$ws.Range("J259").Value = "Statement-non-opinion"
$ws.Range("I262").Value = "sv"
$ws.Range("J262").Value = "Statement-opinion"
$ws.Range("I267").Value = "sd"
$ws.Range("J267").Value = "Statement-non-opinion"
$ws.Range("I279").Value = "%"
$ws.Range("J279").Value = "Uninterpretable"
$ws.Range("I285").Value = "%"
$ws.Range("J285").Value = "Uninterpretable"
$ws.Range("I289").Value = "sv"
$ws.Range("J289").Value = "Statement-opinion"
$ws.Range("I290").Value = "sv"
$ws.Range("J290").Value = "Statement-opinion"
$ws.Range("I291").Value = "%"
$ws.Range("J291").Value = "Uninterpretable"
$ws.Range("I303").Value = "sd"
$ws.Range("J303").Value = "Statement-non-opinion"
$ws.Range("I304").Value = "qy"
$ws.Range("J304").Value = "Yes-No-Question"
$ws.Range("I313").Value = "aa"
$ws.Range("J313").Value = "Agree/Accept"
$ws.Range("I317").Value = "sv"
$ws.Range("J317").Value = "Statement-opinion"
$ws.Range("I318").Value = "sd"
$ws.Range("J318").Value = "Statement-non-opinion"
$ws.Range("I321").Value = "sd"
$ws.Range("J321").Value = "Statement-non-opinion"
$ws.Range("I326").Value = "b"
$ws.Range("J326").Value = "Acknowledge (Backchannel)"
$ws.Range("I335").Value = "aa"
$ws.Range("J335").Value = "Agree/Accept"
$ws.Range("I342").Value = "aa"
$ws.Range("J342").Value = "Agree/Accept"
$ws.Range("I343").Value = "b"
$ws.Range("J343").Value = "Acknowledge (Backchannel)"
$ws.Range("I358").Value = "%"
$ws.Range("J358").Value = "Uninterpretable"
$ws.Range("I377").Value = "sd"
$ws.Range("J377").Value = "Statement-non-opinion"
$ws.Range("I404").Value = "b"
$ws.Range("J404").Value = "Acknowledge (Backchannel)"
